$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Candidate table's field list lives in column C (rows 16-28), with
# column D holding an (empty) adjacent value cell. A new field, "age", is
# being inserted right after "disabilityType" (row 15) and before "dob"
# (row 16) -- i.e. a single-cell "Insert, Shift Cells Down" on column C
# only. Column D (and the other tables in columns A/E/G) are untouched.
#
# Shift column C values down by one row, working from the bottom up so
# nothing is clobbered before it's read. Row 28 doesn't have any data
# yet, so it receives whatever was in row 27 (the last populated row).
for ($r = 28; $r -ge 16; $r--) {
    $srcVal = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r + 1, 3).Value2 = $srcVal
}

# The last entry (previously "batchId-FK" at C28, styled in red) needs to
# keep looking like the other red "- FK" markers once it lands on C29.
$ws.Cells.Item(29, 3).Font.Color = $ws.Cells.Item(24, 1).Font.Color
$ws.Cells.Item(29, 3).Borders.LineStyle = $ws.Cells.Item(28, 3).Borders.LineStyle

# Now insert the new field in the gap we just opened up.
$ws.Range("C16").Value2 = "age"

# Match the author's final selection/cursor position.
$ws.Range("C16").Select()
